$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy style/format from the adjacent header cell (AC1)
# so the new headers share the same bold/border/alignment style, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2..42): Wins=95, Losses=67, Ties=0
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 95
    $ws.Cells.Item($row, 31).Value = 67
    $ws.Cells.Item($row, 32).Value = 0
}
